$d = $word.ActiveDocument

# Locate the paragraph that contains the sentence about version control / PUN2.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*for version control, and PUN2 for online functionality*") {
        $targetPara = $d.Paragraphs($i)
        break
    }
}

$r = $targetPara.Range

# The paragraph is built from 3 runs: "Unity based game. Using " / "Git" /
# " for version control, and PUN2 for online functionality.". The last run
# must become two runs:
#   " for version control, PUN2 for online functionality"
#   " and TMP for better looking UI elements."
# Any direct text insert/delete in this paragraph causes adjacent runs that
# share identical formatting to coalesce, which would also merge the
# untouched "Git" run into its neighbour. To avoid that we temporarily give
# "Git" a distinguishing (no-op) format delta for the duration of the edit.

$fullText = $r.Text
$gitIdx = $fullText.IndexOf("Git")
$gitStart = $r.Start + $gitIdx
$gitEnd = $gitStart + 3
$gitRange = $d.Range($gitStart, $gitEnd)
$gitRange.Font.Bold = 1

# Build the replacement tail (" and TMP for better looking UI elements.")
# with the correct run formatting by copying a FormattedText sample of the
# target run to a scratch spot at the very end of the document (assigning
# FormattedText does not trigger the run-coalescing pass), editing its text
# there, then copying the finished FormattedText back next to our target
# and removing the scratch copy again.
$curText = $targetPara.Range.Text
$sampleIdx = $curText.IndexOf(" for version control")
$sampleStart = $r.Start + $sampleIdx
$sample = $d.Range($sampleStart, $sampleStart + 1)
$sampleFormatted = $sample.FormattedText

$scratchPos = $d.Content.End - 1
$d.Range($scratchPos, $scratchPos).FormattedText = $sampleFormatted
$scratchRange = $d.Range($scratchPos, $d.Content.End - 1)
$scratchRange.Text = " and TMP for better looking UI elements."

# Remove "and " before "PUN2".
$curText = $targetPara.Range.Text
$andIdx = $curText.IndexOf("and PUN2")
$andStart = $r.Start + $andIdx
$andEnd = $andStart + "and ".Length
$d.Range($andStart, $andEnd).Text = ""

# Remove the trailing period after "functionality".
$curText = $targetPara.Range.Text
$periodIdx = $curText.IndexOf("functionality.")
$periodPos = $r.Start + $periodIdx + "functionality".Length
$d.Range($periodPos, $periodPos + 1).Text = ""

# Re-locate the scratch text (its absolute position shifted because of the
# edits above) and copy its FormattedText - carrying the run's original
# formatting - onto the end of the target paragraph as a brand new run.
$needle = " and TMP for better looking UI elements."
$wholeText = $d.Content.Text
$scratchIdxNow = $wholeText.LastIndexOf($needle)
$scratchStartNow = $d.Content.Start + $scratchIdxNow
$scratchEndNow = $scratchStartNow + $needle.Length
$finishedFormatted = $d.Range($scratchStartNow, $scratchEndNow).FormattedText

$paraEnd = $targetPara.Range.End - 1
$d.Range($paraEnd, $paraEnd).FormattedText = $finishedFormatted

# Remove the scratch copy (re-locate once more: inserting the new run above
# shifted everything after it, including the scratch text, forward).
$wholeText2 = $d.Content.Text
$scratchIdx2 = $wholeText2.LastIndexOf($needle)
$scratchStart2 = $d.Content.Start + $scratchIdx2
$scratchEnd2 = $scratchStart2 + $needle.Length
$d.Range($scratchStart2, $scratchEnd2).Text = ""

# Restore "Git"'s original (non-bold) formatting.
$d.Range($gitStart, $gitEnd).Font.Bold = 0

Write-Output $targetPara.Range.Text
